$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("大智投资组合")
for ($r = 2; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 5).Value = "202509211559"
}

$ws2 = $wb.Worksheets.Item("大成投资组合")
for ($r = 2; $r -le 11; $r++) {
    $ws2.Cells.Item($r, 5).Value = "202509211559"
}

$ws3 = $wb.Worksheets.Item("我的投资组合")
for ($r = 2; $r -le 13; $r++) {
    $ws3.Cells.Item($r, 7).Value = "202509211559"
}
